# Template update: rename/translate header columns (Spanish/legacy -> English
# naming), drop the two unused blank sheets, and give the renamed indicator
# headers their own date-ish number-format style (yyyy/mm/dd).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Hoja1")

# --- Rename header cells E1:AD1 (columns A1:D1 keep their existing text) ---
$ws.Range("E1").Value  = "Age_Group"
$ws.Range("F1").Value  = "ILINumFem"
$ws.Range("G1").Value  = "ILINumMale"
$ws.Range("H1").Value  = "ILINumST"
$ws.Range("I1").Value  = "ILINumEmerST"
$ws.Range("J1").Value  = "ILIDenoFem"
$ws.Range("K1").Value  = "ILIDenoMale"
$ws.Range("L1").Value  = "ILIDenoST"
$ws.Range("M1").Value  = "HospFem"
$ws.Range("N1").Value  = "HospMale"
$ws.Range("O1").Value  = "HospST"
$ws.Range("P1").Value  = "ICUFem"
$ws.Range("Q1").Value  = "ICUMale"
$ws.Range("R1").Value  = "ICUST"
$ws.Range("S1").Value  = "DeathsFem"
$ws.Range("T1").Value  = "DeathsMale"
$ws.Range("U1").Value  = "DeathsST"
$ws.Range("V1").Value  = "PneuFem"
$ws.Range("W1").Value  = "PneuMale"
$ws.Range("X1").Value  = "PneuST"
$ws.Range("Y1").Value  = "CCSARIFem"
$ws.Range("Z1").Value  = "CCSARIMale"
$ws.Range("AA1").Value = "CCSARIST"
$ws.Range("AB1").Value = "VentFem"
$ws.Range("AC1").Value = "VentMale"
$ws.Range("AD1").Value = "VentST"

# --- Give the same E1:AD1 header cells their new number format so they pick
#     up the dedicated bold/fill style (mirrors the style split seen on
#     M1, O1, Y1, AA1, AB1, AD1 in the target workbook). ---
$ws.Range("E1:AD1").NumberFormat = "yyyy/mm/dd"

# --- Drop the two unused blank worksheets ---
$wb.Worksheets.Item("Hoja3").Delete() | Out-Null
$wb.Worksheets.Item("Hoja2").Delete() | Out-Null
